# Updated cryptos list on Wed May 29 05:49:16 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns for
# the crypto-ranking rows on the active sheet. Values that would otherwise
# be auto-detected as numbers (losing the exact literal text, e.g. a
# trailing zero like "603.00" -> 603) are written with a leading apostrophe
# so Excel stores them verbatim as text, matching the source data feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.717.89"
$ws.Range("E2").Value = "  +1.15%  "

$ws.Range("D3").Value = "3.876.61"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'603.00"
$ws.Range("E5").Value = "  +0.71%  "

$ws.Range("D6").Value = "'171.71"
$ws.Range("E6").Value = "  +3.13%  "

$ws.Range("D7").Value = "3.876.87"
$ws.Range("E7").Value = "  +0.62%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").Value = "  +1.01%  "

$ws.Range("D11").Value = "'6.52"
$ws.Range("E11").Value = "  +3.40%  "

$ws.Range("E12").Value = "  +1.33%  "

$ws.Range("D13").Value = "'0.0000287"
$ws.Range("E13").Value = "  +15.58%  "

$ws.Range("D14").Value = "'37.22"
$ws.Range("E14").Value = "  +0.81%  "

$ws.Range("D15").Value = "4.525.25"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").Value = "3.864.91"
$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("D17").Value = "68.746.42"
$ws.Range("E17").Value = "  +1.14%  "

$ws.Range("E18").Value = "  +1.10%  "

$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("E21").Value = "  +2.06%  "

$ws.Range("D22").Value = "'472.41"
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("D23").Value = "'0.735"
$ws.Range("E23").Value = "  +0.91%  "

$ws.Range("E24").Value = "  +1.28%  "

$ws.Range("E25").Value = "  +0.98%  "

$ws.Range("E26").Value = "  +2.51%  "

$ws.Range("D27").Value = "'12.29"
$ws.Range("E27").Value = "  +0.94%  "

$ws.Range("D28").Value = "'10.49"
$ws.Range("E28").Value = "  +4.66%  "

$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").Value = "4.026.21"
$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("E32").Value = "  +1.42%  "

$ws.Range("E33").Value = "  +0.80%  "

$ws.Range("D34").Value = "'31.40"
$ws.Range("E34").Value = "  +1.36%  "

$ws.Range("D35").Value = "'9.44"
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").Value = "3.840.77"
$ws.Range("E36").Value = "  +0.22%  "

$ws.Range("D37").Value = "'3.96"
$ws.Range("E37").Value = "  +21.94%  "

$ws.Range("E38").Value = "  +1.22%  "

$ws.Range("D39").Value = "'6.03"
$ws.Range("E39").Value = "  +2.22%  "

$ws.Range("E40").Value = "  +0.50%  "

$ws.Range("E41").Value = "  +0.74%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("E43").Value = "  +2.75%  "

$ws.Range("D44").Value = "'0.000306"
$ws.Range("E44").Value = "  +14.57%  "

$ws.Range("D45").Value = "'2.00"
$ws.Range("E45").Value = "  +1.04%  "

$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").Value = "'8.78"
$ws.Range("E47").Value = "  +2.80%  "

$ws.Range("D48").Value = "'420.64"
$ws.Range("E48").Value = "  -1.38%  "

$ws.Range("D49").Value = "'46.71"
$ws.Range("E49").Value = "  -1.07%  "

$ws.Range("D50").Value = "'142.74"
$ws.Range("E50").Value = "  -0.69%  "

$ws.Range("D51").Value = "'0.0361"
$ws.Range("E51").Value = "  +1.74%  "
